# Saldo_guide.xlsx update
# 1) Shift every "Dt. Referencia" (column G) date forward by one day
#    (45398 -> 45399, i.e. 2024-04-16 -> 2024-04-17) for all data rows.
# 2) Correct the "Saldo Previsto" (D), "Vl. Projetado" (E) and
#    "Vl. Total" (H) figures on a handful of rows.
# 3) Replace the active selection with a "select-all" range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bump every date in column G (rows 2-310) by one day ------------
for ($r = 2; $r -le 310; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = $cell.Value2 + 1
}

# --- 2) Row-specific corrections to D / E / H --------------------------
$ws.Range("D55").Value = 6.23
$ws.Range("E55").Value = 3492.68
$ws.Range("H55").Value = 3498.91

$ws.Range("D61").Value = 119.62
$ws.Range("E61").Value = 11656.59
$ws.Range("H61").Value = 11776.21

$ws.Range("D72").Value = 414.62
$ws.Range("H72").Value = 414.62

$ws.Range("D73").Value = 10687.07
$ws.Range("H73").Value = 10687.07

$ws.Range("D115").Value = 182.05
$ws.Range("E115").Value = 14446.73
$ws.Range("H115").Value = 14628.78

$ws.Range("D120").Value = 299.82
$ws.Range("E120").Value = 38616.22
$ws.Range("H120").Value = 38916.04

$ws.Range("D121").Value = 962.21
$ws.Range("H121").Value = 962.21

$ws.Range("D125").Value = 757.57
$ws.Range("E125").Value = 46275.98
$ws.Range("H125").Value = 47033.55

$ws.Range("D135").Value = 234.67
$ws.Range("H135").Value = 234.67

$ws.Range("D137").Value = 407.6
$ws.Range("H137").Value = 407.6

$ws.Range("E151").Value = 7055.41
$ws.Range("H151").Value = 13684.92

$ws.Range("D285").Value = 7800
$ws.Range("H285").Value = 7800

$ws.Range("D290").Value = 57922.4
$ws.Range("E290").Value = 21333.17
$ws.Range("H290").Value = 79255.57

# --- 3) Selection becomes "select all" (Ctrl+A) -------------------------
$ws.Cells.Select()
